# Normalize column naming to snake_case variable names and Title-Case
# the Spanish location-name prepositions ("de", "del", "el", "la", "las",
# "los", "y") throughout the municipality/state columns. Also trims the
# trailing metadata/footnote rows that shouldn't be part of the clean
# tabular data, and the single "TOTAL" label is normalized to "Total".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header columns (row 1) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- 2. Title-case the Spanish prepositions inside the state/municipality text ---
$preps = @('de', 'del', 'el', 'la', 'las', 'los', 'y')

function Convert-PrepTitleCase($text) {
    $words = $text.Split(' ')
    for ($i = 0; $i -lt $words.Length; $i++) {
        $w = $words[$i]
        if ($w.Length -gt 0 -and ($preps -contains $w.ToLower())) {
            $words[$i] = $w.Substring(0, 1).ToUpper() + $w.Substring(1).ToLower()
        }
    }
    return [string]::Join(' ', $words)
}

$lastDataRow = 1036

for ($r = 2; $r -le $lastDataRow; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $valA = $cellA.Value()
    if ($valA -ne $null -and $valA -is [string] -and $valA.Length -gt 0) {
        $newA = Convert-PrepTitleCase $valA
        $cellA.Value = $newA
    }

    $cellB = $ws.Cells.Item($r, 2)
    $valB = $cellB.Value()
    if ($valB -ne $null -and $valB -is [string] -and $valB.Length -gt 0) {
        $newB = Convert-PrepTitleCase $valB
        $cellB.Value = $newB
    }
}

# --- 3. Normalize the grand-total label ---
$ws.Range("A1036").Value = "Total"

# --- 4. Remove the trailing footnote/metadata rows (1038-1042) ---
$ws.Range("A1038:D1042").EntireRow.Delete()
